# Add two new columns, I (I0) and J (IF), to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---------------------------------------------
# Copy the formatting of the existing header cell H1 ("IP") onto the two
# new header cells so they pick up the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data values (rows 2-62) -------------------------------------------
$I = @(10,6,4,5,6,3,7,3,1,7,8,7,5,7,6,9,8,5,5,6,4,8,7,6,5,7,7,5,7,6,4,4,6,8,6,3,4,7,7,8,7,7,9,7,6,8,9,8,8,3,7,7,7,4,9,5,9,8,1,1,1)
$J = @(10,6,6,6,7,4,8,4,2,7,8,8,5,7,6,9,8,6,6,6,6,8,8,6,6,7,8,6,8,7,5,6,7,9,7,4,6,9,8,8,8,9,9,8,8,9,9,8,8,6,8,8,9,6,9,7,9,8,4,3,2)

for ($i = 0; $i -lt $I.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
}
